# Update Allure test reports by replacing old results with new data:
# - Column G ("act") flips from "Pass" to "Fail" for a handful of rows
#   where the actual result regressed.
# - Column H ("executed") flips from TRUE to FALSE for (almost) every
#   data row, except row 15 which keeps its original executed state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G ("act") : Pass -> Fail ---------------------------------
$ws.Range("G2").Value  = "Fail"
$ws.Range("G12").Value = "Fail"
$ws.Range("G13").Value = "Fail"
$ws.Range("G14").Value = "Fail"

# --- Column H ("executed") : TRUE -> FALSE ---------------------------
# Rows 2-14 and 16-57 (row 15 stays TRUE / unchanged).
$ws.Range("H2:H14").Value  = $false
$ws.Range("H16:H57").Value = $false
